$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New bibliography rows 99-105 (chapter 3 & 4 remake).
# Columns: A=Paper Title, B=Journal, C=Year, D=Author(s), F=Read, G=Point of
# the paper, I=Zotero, J=Available, K=Cited.
# F values: "yes" (style like F9), "no" (style like F98), "~" (style like F97)
# ---------------------------------------------------------------------------

$rows = @(
  @{ Row=99;  A="Integrating animal temperament within ecology and evolution"; B="Biological Reviews"; C=2007; D="Réale et al."; F="approx"; G="THE BIG FIVES (personality)" },
  @{ Row=100; A="Animal personalities: consequences for ecology and evolution"; B="Trends Ecol. Evol."; C=2012; D="Wolf & Weissing"; F="approx"; G="Impact of personality on eco-evo (e.g., pop dynamic) REVIEW" },
  @{ Row=101; A="The development of animal personality: relevance, concepts and perspectives"; B="Biological Reviews"; C=2010; D="Stamps & Groothuis"; F="no"; G="impact of age on personnality" },
  @{ Row=102; A="Animal social networks"; B="Oxford University Press"; C=2015; D="Krause et al."; F="approx"; G="Animal social networks" },
  @{ Row=103; A="Constructing, conducting and interpreting animal social network analysis"; B="Journal of Animal Ecology"; C=2015; D="Farine and Whitehead"; F="yes"; G="SNA" },
  @{ Row=104; A="The Evolution of Bet Hedging in Response to Local Ecological Conditions"; B="The american naturalist"; C=2014; D="Rajon, Desouhant, Chevalier, Débias, Menu"; F="no"; G="bet-hedging" },
  @{ Row=105; A="Environmental unpredicatbility and offspring size: conservative versus diversified bet-hedging"; B="Evolutionary Ecology Research"; C=2004; D="Einum and flemming"; F="approx"; G="Conservative bet-hedging => bigger to buffer unpredictability" }
)

# Source rows/cells whose formatting we clone for each column so the new
# rows look exactly like the rest of the table.
$fmtRow = 98
$fmtCols = @("A","B","C","D","G","I","J","K")
$fmtFYes = "F9"
$fmtFNo = "F98"
$fmtFApprox = "F97"

foreach ($r in $rows) {
  $rowNum = $r.Row

  foreach ($col in $fmtCols) {
    $ws.Range($col + $fmtRow).Copy()
    $ws.Range($col + $rowNum).PasteSpecial(-4122)
  }

  if ($r.F -eq "yes") {
    $ws.Range($fmtFYes).Copy()
  } elseif ($r.F -eq "no") {
    $ws.Range($fmtFNo).Copy()
  } else {
    $ws.Range($fmtFApprox).Copy()
  }
  $ws.Range("F" + $rowNum).PasteSpecial(-4122)

  $ws.Range("A" + $rowNum).Value = $r.A
  $ws.Range("B" + $rowNum).Value = $r.B
  $ws.Range("C" + $rowNum).Value = $r.C
  $ws.Range("D" + $rowNum).Value = $r.D
  $ws.Range("G" + $rowNum).Value = $r.G
  $ws.Range("I" + $rowNum).Value = "yes"
  $ws.Range("J" + $rowNum).Value = "yes"
  $ws.Range("K" + $rowNum).Value = "yes"

  if ($r.F -eq "yes") {
    $ws.Range("F" + $rowNum).Value = "yes"
  } elseif ($r.F -eq "no") {
    $ws.Range("F" + $rowNum).Value = "no"
  } else {
    $ws.Range("F" + $rowNum).Value = [char]0x2248
  }
}

$excel.CutCopyMode = $false

# Recalculate the summary COUNTIF formulas (N19:N26) now that new rows exist.
$excel.CalculateFull()

# Move the cursor/selection to match the author's last position.
$ws.Activate()
$ws.Range("K105").Select()
